# Fixing serial connection bug
# Update the timestamps for rows 65-70, replace row 71 with the former
# row-92 reading (updated timestamp), and drop the now-obsolete rows
# 72-109 that used to follow it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(65, 1).Value = "21-05-25 08:44:03"
$ws.Cells.Item(66, 1).Value = "21-05-25 08:45:51"
$ws.Cells.Item(67, 1).Value = "21-05-25 12:44:57"
$ws.Cells.Item(68, 1).Value = "21-05-25 13:03:49"
$ws.Cells.Item(69, 1).Value = "21-05-25 13:04:12"
$ws.Cells.Item(70, 1).Value = "21-05-25 13:26:47"

# Row 71 takes the place of the old row 92, now with an updated
# timestamp; rows 72-109 are removed entirely.
$ws.Cells.Item(71, 1).Value = "21-05-25 13:27:07"

$ws.Range("A72:E109").EntireRow.Delete()
